$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '23.894.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.08%  '
$ws.Range("D3").Value = "'" + '1.649.89'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = "'" + '310.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").Value = "'" + '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("E7").Value = '  -1.97%  '
$ws.Range("D8").Value = "'" + '0.3816'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.30%  '
$ws.Range("D9").Value = "'" + '51.91'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.12%  '
$ws.Range("D10").Value = "'" + '1.347'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.81%  '
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").Value = "'" + '0.08458'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.51%  '
$ws.Range("D13").Value = "'" + '23.91'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.12%  '
$ws.Range("D14").Value = "'" + '7.060'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.50%  '
$ws.Range("D15").Value = "'" + '8.022'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.68%  '
$ws.Range("D16").Value = "'" + '0.00001312'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.63%  '
$ws.Range("D17").Value = "'" + '1.652.28'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").Value = "'" + '94.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").Value = "'" + '0.07016'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("D20").Value = "'" + '19.66'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.48%  '
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("D23").Value = "'" + '13.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").Value = "'" + '23.895.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.04%  '
$ws.Range("D25").Value = "'" + '2.456'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.42%  '
$ws.Range("E26").Value = '  -3.61%  '
$ws.Range("D27").Value = "'" + '22.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("D28").Value = "'" + '153.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.45%  '
$ws.Range("D29").Value = "'" + '5.408'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("D30").Value = "'" + '138.04'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.29%  '
$ws.Range("D31").Value = "'" + '7.820'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.29%  '
$ws.Range("D32").Value = "'" + '2.502'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.23%  '
$ws.Range("D33").Value = "'" + '1.831.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.79%  '
$ws.Range("D34").Value = "'" + '1.017'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.63%  '
$ws.Range("D35").Value = "'" + '0.08187'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = "'" + '0.02936'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.54%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = "'" + '6.710'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.67%  '
$ws.Range("D38").Value = "'" + '0.2679'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.98%  '
$ws.Range("D39").Value = "'" + '10.71'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.46%  '
$ws.Range("D40").Value = "'" + '0.09126'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.08%  '
$ws.Range("D41").Value = "'" + '0.7579'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.11%  '
$ws.Range("D42").Value = "'" + '13.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.87%  '
$ws.Range("D43").Value = "'" + '1.421'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.50%  '
$ws.Range("D44").Value = "'" + '16.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.21%  '
$ws.Range("D45").Value = "'" + '0.6945'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.55%  '
$ws.Range("D46").Value = "'" + '2.459'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.19%  '
$ws.Range("D47").Value = "'" + '4.091'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.99%  '
$ws.Range("D48").Value = "'" + '1.002'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("D49").Value = "'" + '0.08282'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.95%  '
$ws.Range("D50").Value = "'" + '134.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'" + '1.232'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.94%  '
